$wb = $excel.ActiveWorkbook

# 1) Rename sheets so the country-code prefix is followed by a space:
#    "XXXNOCB" -> "XXX NOCB", "XXXRolling statistics" -> "XXX Rolling statistics",
#    "XXXInterpolate" -> "XXX Interpolate"
$countries = @("CHN", "GBR", "HKG", "IND", "IRN", "ITA", "RUS", "USA", "VNM", "ZAF")
$suffixes = @("NOCB", "Rolling statistics", "Interpolate")

foreach ($country in $countries) {
    foreach ($suffix in $suffixes) {
        $oldName = "$country$suffix"
        $newName = "$country $suffix"
        $ws = $wb.Worksheets.Item($oldName)
        $ws.Name = $newName
    }
}

# 2) Fix the "Rolling statistics" sheets' column A: rows 3-43 hold fractional
#    pseudo-year values that should be the plain integer year (row 3 -> 1981,
#    row 4 -> 1982, ..., row 43 -> 2021).
foreach ($country in $countries) {
    $ws = $wb.Worksheets.Item("$country Rolling statistics")
    for ($row = 3; $row -le 43; $row++) {
        $year = $row + 1978
        $ws.Cells.Item($row, 1).Value = $year
    }
}
